# Automatic update of files.
# Rotates the species-record rows 2..4:
#   row2 <- old row4 data, row3 <- old row2 data, row4 <- old row3 data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, $val) {
    # Force text storage (even for numeric-looking strings like "1"/"2"),
    # then restore the default "Normal" style so no stray number format sticks.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Row 2 (becomes old row 4's data) ---
$ws.Range("A2").Value = 68176169
$ws.Range("B2").Value = 57587
Set-TextCell $ws.Range("C2") "Behöver inte valideras"
$ws.Range("E2").Value = 100141
Set-TextCell $ws.Range("F2") "Större vattensalamander"
Set-TextCell $ws.Range("G2") "Triturus cristatus"
Set-TextCell $ws.Range("H2") "(Laurenti, 1768)"
Set-TextCell $ws.Range("I2") ""
Set-TextCell $ws.Range("J2") ""
Set-TextCell $ws.Range("L2") ""
Set-TextCell $ws.Range("P2") "Sjövik, Sm"
$ws.Range("Q2").Value = 503526.6896539551
$ws.Range("R2").Value = 6429839.084042171
Set-TextCell $ws.Range("Y2") "2005-01-01"
Set-TextCell $ws.Range("AA2") "2005-12-31"
Set-TextCell $ws.Range("AX2") "Anna Isaksson"

# --- Row 3 (becomes old row 2's data) ---
$ws.Range("A3").Value = 68175904
$ws.Range("B3").Value = 57585
$ws.Range("E3").Value = 208242
Set-TextCell $ws.Range("F3") "Mindre vattensalamander"
Set-TextCell $ws.Range("G3") "Lissotriton vulgaris"
Set-TextCell $ws.Range("H3") "(Linnaeus, 1758)"
Set-TextCell $ws.Range("I3") "2"

# --- Row 4 (becomes old row 3's data) ---
$ws.Range("A4").Value = 68175906
Set-TextCell $ws.Range("C4") "Ovaliderad"
Set-TextCell $ws.Range("I4") "1"
Set-TextCell $ws.Range("J4") "ex."
Set-TextCell $ws.Range("L4") "hona"
Set-TextCell $ws.Range("P4") "Sjövik, Sommen, Sm"
$ws.Range("Q4").Value = 503498.5757228022
$ws.Range("R4").Value = 6429815.746484536
Set-TextCell $ws.Range("Y4") "1999-04-28"
Set-TextCell $ws.Range("AA4") "1999-04-28"
Set-TextCell $ws.Range("AX4") "Josefine Gustafsson"
